$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.446.45"
$ws.Range("E2").Value = "  -2.87%  "

$ws.Range("D3").Value = "1.806.70"
$ws.Range("E3").Value = "  -2.51%  "

$ws.Range("E4").Value = "  +0.77%  "

$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4555"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.61%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3660"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07132"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8764"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.28%  "

$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.61%  "

$ws.Range("D13").Value = "1.771.06"
$ws.Range("E13").Value = "  -9.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.269"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.367"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "85.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.79%  "

$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008564"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.93%  "

$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D20").Value = "26.502.41"
$ws.Range("E20").Value = "  -2.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.986"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.983"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.042"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "112.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.833"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08664"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.042"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.93%  "

$ws.Range("E32").Value = "  -4.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.463"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.112"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.05%  "

$ws.Range("E35").Value = "  +0.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.509"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.078"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01931"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05106"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.893"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.942"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5018"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1561"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.120"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.007"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4611"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.984"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.587"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06005"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.17%  "
